$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.450.54'
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").Value = '2.664.54'
$ws.Range("E3").Value = '  +1.72%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.45'
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.33'
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.175'
$ws.Range("E9").Value = '  +6.01%  '

$ws.Range("D10").Value = '2.664.17'
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("E11").Value = '  +1.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.360'
$ws.Range("E12").Value = '  +4.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.06'
$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.151.49'
$ws.Range("E14").Value = '  +1.69%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000191'
$ws.Range("E15").Value = '  +4.19%  '

$ws.Range("D16").Value = '72.345.83'
$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.66'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").Value = '2.667.95'
$ws.Range("E18").Value = '  +2.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.06'
$ws.Range("E19").Value = '  +4.97%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.00'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.30'
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.21'
$ws.Range("E22").Value = '  +1.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.08'
$ws.Range("E23").Value = '  +12.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.84'
$ws.Range("E24").Value = '  +1.28%  '

$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.42'
$ws.Range("E25").Value = '  -0.35%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +5.05%  '

$ws.Range("D28").Value = '2.807.14'
$ws.Range("E28").Value = '  +1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = '0.0₃0959'
$ws.Range("E30").Value = '  +0.68%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.23'
$ws.Range("E31").Value = '  +2.91%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '523.80'
$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.32'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.83'
$ws.Range("E36").Value = '  +0.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.55'
$ws.Range("E37").Value = '  +2.31%  '

$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.12'
$ws.Range("E38").Value = '  +0.72%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.40'
$ws.Range("E39").Value = '  +2.18%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.112'
$ws.Range("E40").Value = '  -6.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.85'
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.12'
$ws.Range("E42").Value = '  +2.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.63'
$ws.Range("E43").Value = '  +1.55%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.31'
$ws.Range("E46").Value = '  -2.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.68'
$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.74'
$ws.Range("E48").Value = '  +2.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.550'
$ws.Range("E49").Value = '  +3.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").Value = '  +3.13%  '

$ws.Range("D51").Value = '0.0₆0262'
$ws.Range("E51").Value = '  -1.81%  '
